$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 data. A3 looks like a date ("2025-10-18"); force it to stay as
# plain text (matching the inlineStr/text cell in the target) instead of
# letting Excel auto-convert it into a date serial number, then restore
# the default "Normal" style so no extra number-format style is left on
# the cell.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-10-18"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "YYY"
$ws.Range("C3").Value = "123ABX007"
$ws.Range("D3").Value = "Karapakkam"
